$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 23.15000000000018
$ws.Range("H2").Value = [double]"2.766910964797898e-16"
$ws.Range("K2").Value = 46.85925603551192
$ws.Range("L2").Value = "[39.703209297831535, 54.0153027731923]"
$ws.Range("O2").Value = 1.553500271144502
$ws.Range("P2").Value = "[1.3899739268135018, 1.7170266154755023]"
$ws.Range("S2").Value = 59.1586710533139
$ws.Range("T2").Value = "[54.27185372466471, 64.04548838196308]"
$ws.Range("W2").Value = 17.42622622622636
$ws.Range("X2").Value = 16.82372372372386
$ws.Range("Y2").Value = 18.02872872872887

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.21000000000003
$ws.Range("H3").Value = [double]"2.766910964797898e-16"
$ws.Range("K3").Value = 48.85352682776533
$ws.Range("L3").Value = "[34.9823716482893, 62.724682007241356]"
$ws.Range("M3").Value = [double]"5.838507455280251e-11"
$ws.Range("N3").Value = [double]"5.838507455280251e-11"
$ws.Range("O3").Value = -0.9308422677303092
$ws.Range("P3").Value = "[-1.2075791581366175, -0.6541053773240009]"
$ws.Range("Q3").Value = [double]"3.32595728735896e-10"
$ws.Range("R3").Value = [double]"3.32595728735896e-10"
$ws.Range("S3").Value = 63.32893697685498
$ws.Range("T3").Value = "[56.13475630477417, 70.5231176489358]"
$ws.Range("W3").Value = 3.290370370370375
$ws.Range("X3").Value = 2.312152152152155
$ws.Range("Y3").Value = 4.268588588588596
